$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "код идентификатор PLU"
$ws.Range("B1").Value = "наименование анализируемых позиций"
$ws.Range("C1").Value = "данные по анализируемому критерию (продажи/оборот/прибыль)"

$hdr = $ws.Range("A1:C1")
$hdr.Font.Size = 9
$hdr.Font.Color = 0
$hdr.HorizontalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Rows(1).RowHeight = 36.6

# ---- Data rows (rows 2-6) ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Товар 1"
$ws.Range("C2").Value = 100

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Товар 2"
$ws.Range("C3").Value = 50

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Товар 3"
$ws.Range("C4").Value = ""

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = 20

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Товар 5"
$ws.Range("C6").Value = 5

$data = $ws.Range("A2:C6")
$data.Font.Name = "Arial"
$data.Font.Size = 10

# ---- Column widths ----
$ws.Columns("A:B").ColumnWidth = 33.5
$ws.Columns("C").ColumnWidth = 56.05

# ---- Selection shown in the saved view ----
$ws.Range("B12").Select()
